$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new "Baseline + Glove Word vectors" result row
$ws.Range("A6").Value = "Baseline + Glove Word vectors"
$ws.Range("B6").Value = "InceptionV3 (2048 units)"
$ws.Range("C6").Value = "LSTM (2048 units)"
$ws.Range("E6").Value = 0.3592
$ws.Range("F6").Value = 0.3968
$ws.Range("G6").Value = 0.5893
$ws.Range("H6").Value = 0.6715

# Row 7: new "Baseline + Glove Word vectors + 2x LSTM" result row
$ws.Range("A7").Value = "Baseline + Glove Word vectors + 2x LSTM"
$ws.Range("B7").Value = "InceptionV3 (2048 units)"
$ws.Range("C7").Value = "LSTM (2048 units)"
$ws.Range("E7").Value = 0.3515
$ws.Range("F7").Value = 0.2503
$ws.Range("G7").Value = 0.4028
$ws.Range("H7").Value = 0.2631

# Distinguish the A7 label with its own (Calibri 11 black) font, matching the
# author's formatting tweak on that row.
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Bold = $false

# Move the active selection, as recorded after the edit.
$ws.Range("E8").Select()
